# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet - this shifts the old "Late" column (N) to O and the old
# "Outstanding" column (P) to Q - then give the newly inserted column the
# same width as its neighbour "In Advance" (M), and finally leave the
# "Repayment Schedule" sheet active with cell R6 selected (this also moves
# the active-tab/tabSelected flags off "NewLoanInput" and onto
# "Repayment Schedule").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

$ws.Activate() | Out-Null
$ws.Range("R6").Select() | Out-Null
